$wb = $excel.ActiveWorkbook

# =====================================================================
# Change: dataset-shape labels printed as Python tuple-repr strings,
# e.g. "(5554, 145)", "(20,)", "()" are rewritten into Excel/VB style
# array-literal notation, e.g. "{ 5554, 145 }", "{ 20 }", "{}".
# =====================================================================

# --- h5showList3 ---
$ws = $wb.Worksheets.Item("h5showList3")
$ws.Range("H4").Value = "{ 21684 }"
$ws.Range("H5").Value = "{ 2470 }"
$ws.Range("H6").Value = "{ 14711 }"
$ws.Range("H7").Value = "{ 23536 }"
$ws.Range("H8").Value = "{ 16954 }"
$ws.Range("H9").Value = "{ 17969 }"

# --- h5showList ---
$ws = $wb.Worksheets.Item("h5showList")
foreach ($c in @("H10","H11","H13","H15","H16","H18","H20","H21","H23","H25","H26","H28","H30","H31","H33","H35","H37","H39","H41","H43","H45","H47","H49","H52","H54","H56","H58","H60","H61","H63","H66","H68","H70","H72","H74","H76","H78","H81")) { $ws.Range($c).Value = "{ 5554, 145 }" }
foreach ($c in @("H12","H14","H17","H19","H22","H24","H27","H29","H32","H34","H36","H38","H40","H42","H44","H46","H48","H50","H53","H55","H57","H59","H62","H64","H65","H67","H69","H71","H73","H75","H77","H79","H82","H83","H84","H85","H87","H88","H89","H90","H91","H92","H93","H94","H95","H96","H97","H98","H99","H100")) { $ws.Range($c).Value = "{ 5554 }" }
$ws.Range("H51").Value = "{ 5554, 1 }"
$ws.Range("H86").Value = "{ 145 }"
foreach ($c in @("H102","H103")) { $ws.Range($c).Value = "{}" }

# --- h5showList1 ---
$ws = $wb.Worksheets.Item("h5showList1")
$ws.Range("H6").Value = "{ 10, 10 }"
$ws.Range("H7").Value = "{ 20 }"
$ws.Range("H11").Value = "{ 10 }"
$ws.Range("H12").Value = "{ 3, 5 }"

# =====================================================================
# Change: record a selection of A2 on h5showList2 and h5showList3 (the
# sheets that didn't already carry a saved selection), and set the
# h5showList3 page orientation to portrait. h5showList3 stays the
# active/selected tab, so its selection is applied last.
# =====================================================================
$ws4 = $wb.Worksheets.Item("h5showList2")
$ws4.Range("A2").Select()

$ws5 = $wb.Worksheets.Item("h5showList3")
$ws5.PageSetup.Orientation = 1
$ws5.Range("A2").Select()
